$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("w","z","u","s","O_2x","O_2y","O_4x","O_4y","A1X","A1Y","B1X","B1Y","A2X","A2Y","B2X","B2Y","G1","V1","W + V1 + U","THETA","PHI","BETA2","SIGMA","PSI","GAMMA2")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Range("Y1").Select()
